$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Ngf","Ngfr","FAPs",1,0.3333333333333333,0.2000683333333333,0.600205,0.01182207550347303,0.01182207550347303,3,1,3.362744666666666,10.088234,0.7488888671136141,0.748888867113614,0.6727787208855555,6.055008487969999,0.008853420730727524,0.008853420730727523),
    @("ECs","Ngf","Ngfr","sCs",1,0.3333333333333333,0.2000683333333333,0.600205,0.01182207550347303,0.01182207550347303,3,1,1.127567333333333,3.382702,0.251111132886386,0.2511111328863859,0.2255905171011111,2.03031465391,0.002968654772745504,0.002968654772745503),
    @("FAPs","Ngf","Ngfr","FAPs",1,0.3333333333333333,0.44525,1.33575,0.02630990637159653,0.02630990637159653,3,1,3.362744666666666,10.088234,0.7488888671136141,0.748888867113614,1.497262062833333,13.4753585655,0.01970319597649018,0.01970319597649018),
    @("FAPs","Ngf","Ngfr","sCs",1,0.3333333333333333,0.44525,1.33575,0.02630990637159653,0.02630990637159653,3,1,1.127567333333333,3.382702,0.251111132886386,0.2511111328863859,0.5020493551666667,4.5184441965,0.00660671039510635,0.006606710395106348),
    @("sCs","Ngf","Ngfr","FAPs",3,1,16.277965,48.833895,0.9618680181249305,0.9618680181249305,3,1,3.362744666666666,10.088234,0.7488888671136141,0.748888867113614,54.73863998793666,492.64775989143,0.7203322504063965,0.7203322504063964),
    @("sCs","Ngf","Ngfr","sCs",3,1,16.277965,48.833895,0.9618680181249305,0.9618680181249305,3,1,1.127567333333333,3.382702,0.251111132886386,0.2511111328863859,18.35450158714333,165.19051428429,0.2415357677185341,0.2415357677185341)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
